# Updated cryptos list values (Price / Volume(1h), and for rows 41-44 also Coin/Link)
# to match the latest scrape, per the authoring diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '43.872.79'
$ws.Range('E2').Value = '  -0.38%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.357.26'
$ws.Range('E3').Value = '  -0.16%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.02%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '239.80'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.02%  '

# Row 6: XRP
$ws.Range('E6').Value = '  -1.74%  '

# Row 7: Solana
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '73.51'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.19%  '

# Row 9: Cardano
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.603'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.73%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  +1.57%  '

# Row 11: OKB
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '60.82'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.33%  '

# Row 12: Avalanche
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '34.01'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.05%  '

# Row 13: TRON
$ws.Range('E13').Value = '  +0.36%  '

# Row 14: Polkadot
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.22'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.01%  '

# Row 15: Chainlink
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '16.18'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.01%  '

# Row 16: Polygon
$ws.Range('E16').Value = '  +0.65%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '2.356.70'
$ws.Range('E17').Value = '  -0.21%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '43.858.68'
$ws.Range('E18').Value = '  -0.14%  '

# Row 19: ShibaInu
$ws.Range('E19').Value = '  +0.63%  '

# Row 20: Litecoin
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '77.78'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.92%  '

# Row 21: Uniswap
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.53'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.60%  '

# Row 22: BitcoinCash
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '252.49'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.81%  '

# Row 23: Dai
$ws.Range('E23').Value = '  -0.03%  '

# Row 24: WEMIXToken
$ws.Range('E24').Value = '  +3.03%  '

# Row 25: ImmutableX
$ws.Range('E25').Value = '  -6.47%  '

# Row 26: PancakeSwap
$ws.Range('E26').Value = '  -0.15%  '

# Row 27: Cosmos
$ws.Range('E27').Value = '  -2.89%  '

# Row 28: Toncoin
$ws.Range('E28').Value = '  +1.57%  '

# Row 29: Monero
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '176.04'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.44%  '

# Row 30: EthereumClassic
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '22.27'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.27%  '

# Row 31: Kaspa
$ws.Range('E31').Value = '  +0.17%  '

# Row 32: Stellar
$ws.Range('E32').Value = '  -2.47%  '

# Row 33: Hedera
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0746'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.57%  '

# Row 34: Filecoin
$ws.Range('E34').Value = '  -3.94%  '

# Row 35: InternetComputer(DFINITY)
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.33'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.42%  '

# Row 36: RenderToken
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.80'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.10%  '

# Row 37: THORChain
$ws.Range('E37').Value = '  +4.48%  '

# Row 38: LidoDAOToken
$ws.Range('E38').Value = '  +1.68%  '

# Row 39: VeChain
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0274'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.47%  '

# Row 40: FTXToken
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.47'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +16.00%  '

# Row 41: MultiversX -> InjectiveProtocol
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '20.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.60%  '

# Row 42: InjectiveProtocol -> MultiversX
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '64.51'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +9.80%  '

# Row 43: Cronos -> FraxShare
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '9.03'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.59%  '

# Row 44: FraxShare -> Cronos
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.106'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -5.90%  '

# Row 45: Algorand
$ws.Range('E45').Value = '  -2.40%  '

# Row 46: BinanceUSD
$ws.Range('E46').Value = '  +0.04%  '

# Row 47: TrustWalletToken
$ws.Range('E47').Value = '  -0.85%  '

# Row 48: NEARProtocol
$ws.Range('E48').Value = '  -2.01%  '

# Row 49: ARBITRUM
$ws.Range('E49').Value = '  -1.98%  '

# Row 50: Aave
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '98.01'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.36%  '

# Row 51: TerraClassic
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.000211'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +13.53%  '
